$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'42.867.75"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "'  -0.57%  "
$ws.Range("E2").Style = "Normal"
$ws.Range("D3").Value = "'2.300.89"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "'  -0.22%  "
$ws.Range("E3").Style = "Normal"
$ws.Range("E4").Value = "'  -0.05%  "
$ws.Range("E4").Style = "Normal"
$ws.Range("D5").Value = "'305.63"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "'  +1.53%  "
$ws.Range("E5").Style = "Normal"
$ws.Range("D6").Value = "'97.14"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "'  -0.60%  "
$ws.Range("E6").Style = "Normal"
$ws.Range("D7").Value = "'0.511"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "'  -1.69%  "
$ws.Range("E7").Style = "Normal"
$ws.Range("E8").Value = "'  -0.07%  "
$ws.Range("E8").Style = "Normal"
$ws.Range("E9").Value = "'  -2.64%  "
$ws.Range("E9").Style = "Normal"
$ws.Range("E10").Value = "'  -0.67%  "
$ws.Range("E10").Style = "Normal"
$ws.Range("E11").Value = "'  +0.37%  "
$ws.Range("E11").Style = "Normal"
$ws.Range("D12").Value = "'18.23"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "'  +1.61%  "
$ws.Range("E12").Style = "Normal"
$ws.Range("E13").Value = "'  +1.06%  "
$ws.Range("E13").Style = "Normal"
$ws.Range("E14").Value = "'  -1.68%  "
$ws.Range("E14").Style = "Normal"
$ws.Range("D15").Value = "'2.661.28"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "'  -0.16%  "
$ws.Range("E15").Style = "Normal"
$ws.Range("D16").Value = "'2.301.82"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "'  +0.39%  "
$ws.Range("E16").Style = "Normal"
$ws.Range("D17").Value = "'0.782"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "'  -0.88%  "
$ws.Range("E17").Style = "Normal"
$ws.Range("D18").Value = "'42.794.72"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "'  -0.49%  "
$ws.Range("E18").Style = "Normal"
$ws.Range("D19").Value = "'13.22"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "'  +0.87%  "
$ws.Range("E19").Style = "Normal"
$ws.Range("E20").Value = "'  -0.50%  "
$ws.Range("E20").Style = "Normal"
$ws.Range("E21").Value = "'  -1.33%  "
$ws.Range("E21").Style = "Normal"
$ws.Range("D22").Value = "'67.57"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "'  -1.14%  "
$ws.Range("E22").Style = "Normal"
$ws.Range("D23").Value = "'236.51"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "'  -0.59%  "
$ws.Range("E23").Style = "Normal"
$ws.Range("D24").Value = "'2.18"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "'  -2.03%  "
$ws.Range("E24").Style = "Normal"
$ws.Range("E25").Value = "'  +2.28%  "
$ws.Range("E25").Style = "Normal"
$ws.Range("E26").Value = "'  +0.86%  "
$ws.Range("E26").Style = "Normal"
$ws.Range("D28").Value = "'25.44"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "'  +0.94%  "
$ws.Range("E28").Style = "Normal"
$ws.Range("D29").Value = "'167.06"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "'  +0.03%  "
$ws.Range("E29").Style = "Normal"
$ws.Range("E30").Value = "'  +1.21%  "
$ws.Range("E30").Style = "Normal"
$ws.Range("E31").Value = "'  -0.84%  "
$ws.Range("E31").Style = "Normal"
$ws.Range("D32").Value = "'33.02"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "'  -0.15%  "
$ws.Range("E32").Style = "Normal"
$ws.Range("E33").Value = "'  +0.05%  "
$ws.Range("E33").Style = "Normal"
$ws.Range("D34").Value = "'4.81"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "'  +1.68%  "
$ws.Range("E34").Style = "Normal"
$ws.Range("E35").Value = "'  -2.57%  "
$ws.Range("E35").Style = "Normal"
$ws.Range("D36").Value = "'17.36"
$ws.Range("D36").Style = "Normal"
$ws.Range("E37").Value = "'  -0.23%  "
$ws.Range("E37").Style = "Normal"
$ws.Range("D38").Value = "'0.0692"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "'  +0.27%  "
$ws.Range("E38").Style = "Normal"
$ws.Range("E39").Value = "'  -0.15%  "
$ws.Range("E39").Style = "Normal"
$ws.Range("E40").Value = "'  -1.74%  "
$ws.Range("E40").Style = "Normal"
$ws.Range("E41").Value = "'  -1.03%  "
$ws.Range("E41").Style = "Normal"
$ws.Range("E42").Value = "'  -1.16%  "
$ws.Range("E42").Style = "Normal"
$ws.Range("D43").Value = "'2.012.58"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "'  +0.11%  "
$ws.Range("E43").Style = "Normal"
$ws.Range("E44").Value = "'  -2.31%  "
$ws.Range("E44").Style = "Normal"
$ws.Range("D45").Value = "'18.19"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "'  +4.10%  "
$ws.Range("E45").Style = "Normal"
$ws.Range("D46").Value = "'10.01"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "'  -2.58%  "
$ws.Range("E46").Style = "Normal"
$ws.Range("E47").Value = "'  -2.99%  "
$ws.Range("E47").Style = "Normal"
$ws.Range("E48").Value = "'  -1.91%  "
$ws.Range("E48").Style = "Normal"
$ws.Range("D49").Value = "'2.94"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "'  +7.73%  "
$ws.Range("E49").Style = "Normal"
$ws.Range("D50").Value = "'53.93"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "'  -1.05%  "
$ws.Range("E50").Style = "Normal"
$ws.Range("D51").Value = "'2.529.65"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "'  -0.30%  "
$ws.Range("E51").Style = "Normal"
